# Generate Report for Handoff
#
# Adds two new file entries to the localization status report:
#   1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f  -> inserted before the existing
#                                            a0dfcb73 row (both "Ready for handoff")
#   bc637f38-e25b-43f1-a14c-3866ee49e0a8  -> appended as a brand new row
#                                            ("Ready for handoff")
#
# Applied identically to the "Overview" sheet and the per-language detail
# sheets ("zh-cn" / "de-de").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Row 5 used to be a0dfcb73's row; it is replaced by the new 1ebd9f8b row,
# a0dfcb73 moves down to row 6, and bc637f38 is appended as row 7.
$overview.Hyperlinks.Add($overview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c8a0c8f2e9a4f6b5d7c1e3a9f0b2d4c6e8a1b3c5/e2e/1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f.md", "", "", "1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f.md")
$overview.Range("B5").Value = "Ready for handoff"
$overview.Range("C5").Value = "Ready for handoff"
$overview.Range("D5").Value = "2016-30-19 06:30:37"

$overview.Hyperlinks.Add($overview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/a0dfcb73-3710-42e5-b5b0-373815c853ed.md", "", "", "a0dfcb73-3710-42e5-b5b0-373815c853ed.md")
$overview.Range("B6").Value = "Ready for handoff"
$overview.Range("C6").Value = "Ready for handoff"
$overview.Range("D6").Value = "2016-28-19 06:28:38"

$overview.Hyperlinks.Add($overview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/f3b6d8a1c5e9f2b4d6a8c0e2f4b6d8a0c2e4f6b8/e2e/bc637f38-e25b-43f1-a14c-3866ee49e0a8.md", "", "", "bc637f38-e25b-43f1-a14c-3866ee49e0a8.md")
$overview.Range("B7").Value = "Ready for handoff"
$overview.Range("C7").Value = "Ready for handoff"
$overview.Range("D7").Value = "2016-30-19 06:30:37"

# ---------------------------------------------------------------------------
# Per-language detail sheets: zh-cn, de-de
# Columns: A Source File Name | B File Extension | C Status |
#          D Latest Handoff File | E Latest Handoff Datetime |
#          F Latest Target File | G Latest Handback File |
#          H Latest Handback DateTime | I Handoff Reason |
#          J Dependency From | K Error Detail
# (rows for files that have not yet been handed off only populate A,B,C,D,E,H,I)
# ---------------------------------------------------------------------------

function Add-LocRow($sheet, $row, $uuid, $mdUrl, $xlfHash, $xlfUrl, $lang, $handoffDatetime) {
    $mdName = "$uuid.md"
    $xlfName = "$uuid.$xlfHash.$lang.xlf"

    $sheet.Hyperlinks.Add($sheet.Range("A$row"), $mdUrl, "", "", $mdName)
    $sheet.Hyperlinks.Add($sheet.Range("B$row"), $mdUrl, "", "", ".md")
    $sheet.Range("C$row").Value = "Ready for handoff"
    $sheet.Hyperlinks.Add($sheet.Range("D$row"), $xlfUrl, "", "", $xlfName)
    $sheet.Range("E$row").Value = $handoffDatetime
    $sheet.Range("H$row").Value = "0001-01-01 00:00:00"
    $sheet.Range("I$row").Value = "Include"
}

$zhcn = $wb.Worksheets.Item("zh-cn")

Add-LocRow $zhcn 5 "1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c8a0c8f2e9a4f6b5d7c1e3a9f0b2d4c6e8a1b3c5/e2e/1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f.md" `
    "640a26003066713c27fbe37a5dbd1591eec56870" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1b2c3d4e5f60718293a4b5c6d7e8f90a1b2c3d4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f.640a26003066713c27fbe37a5dbd1591eec56870.zh-cn.xlf" `
    "zh-cn" `
    "2016-03-19 06:30:32"

Add-LocRow $zhcn 6 "a0dfcb73-3710-42e5-b5b0-373815c853ed" `
    "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/a0dfcb73-3710-42e5-b5b0-373815c853ed.md" `
    "7c23583881e90434debdf5bd12e534d97478fab2" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b40431e5f8e51d6cdae64b193740bad9d014da95/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a0dfcb73-3710-42e5-b5b0-373815c853ed.7c23583881e90434debdf5bd12e534d97478fab2.zh-cn.xlf" `
    "zh-cn" `
    "2016-03-19 06:28:35"

Add-LocRow $zhcn 7 "bc637f38-e25b-43f1-a14c-3866ee49e0a8" `
    "https://github.com/OpenLocalizationTest/oltest/blob/f3b6d8a1c5e9f2b4d6a8c0e2f4b6d8a0c2e4f6b8/e2e/bc637f38-e25b-43f1-a14c-3866ee49e0a8.md" `
    "08ea4fd06e16806ac3e2fb8b5d7db844d55545db" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5d6e7f8a9b0c1d2e3f4a5b6c7d8e9f0a1b2c3d4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bc637f38-e25b-43f1-a14c-3866ee49e0a8.08ea4fd06e16806ac3e2fb8b5d7db844d55545db.zh-cn.xlf" `
    "zh-cn" `
    "2016-03-19 06:30:32"

$dede = $wb.Worksheets.Item("de-de")

Add-LocRow $dede 5 "1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c8a0c8f2e9a4f6b5d7c1e3a9f0b2d4c6e8a1b3c5/e2e/1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f.md" `
    "640a26003066713c27fbe37a5dbd1591eec56870" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d6e7f8a9b0c1d2e3f4a5b6c7d8e9f0a1b2c3d4e5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f.640a26003066713c27fbe37a5dbd1591eec56870.de-de.xlf" `
    "de-de" `
    "2016-03-19 06:30:37"

Add-LocRow $dede 6 "a0dfcb73-3710-42e5-b5b0-373815c853ed" `
    "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/a0dfcb73-3710-42e5-b5b0-373815c853ed.md" `
    "7c23583881e90434debdf5bd12e534d97478fab2" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef0aef0432d93019bea41c1cc46a73929fdaa4fc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a0dfcb73-3710-42e5-b5b0-373815c853ed.7c23583881e90434debdf5bd12e534d97478fab2.de-de.xlf" `
    "de-de" `
    "2016-03-19 06:28:38"

Add-LocRow $dede 7 "bc637f38-e25b-43f1-a14c-3866ee49e0a8" `
    "https://github.com/OpenLocalizationTest/oltest/blob/f3b6d8a1c5e9f2b4d6a8c0e2f4b6d8a0c2e4f6b8/e2e/bc637f38-e25b-43f1-a14c-3866ee49e0a8.md" `
    "08ea4fd06e16806ac3e2fb8b5d7db844d55545db" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f8a9b0c1d2e3f4a5b6c7d8e9f0a1b2c3d4e5f6a7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bc637f38-e25b-43f1-a14c-3866ee49e0a8.08ea4fd06e16806ac3e2fb8b5d7db844d55545db.de-de.xlf" `
    "de-de" `
    "2016-03-19 06:30:37"

$wb.Save()
